$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$urls = @(
    "https://www.facebook.com/thongtinchinhphu/posts/986186283546152?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985610296937084?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985437616954352?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985133396984774?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985123206985793?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985004743664306?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/985110820320365?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/984913873673393?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/984827680348679?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/984231030408344?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/984220663742714?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/983512893813491?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/983100937188020?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/983244920506955?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/982341943930586?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/982678903896890?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981956287302485?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981804130651034?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981741710657276?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981790440652403?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981610670670380?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981385857359528?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981043764060404?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981046457393468?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981039514060829?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/980992030732244?ref=embed_post",
    "https://www.facebook.com/thongtinchinhphu/posts/981154850715962?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/881783077373891?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/881654814053384?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/881313624087503?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/881315360753996?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/881314527420746?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880753314143534?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880752837476915?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880636624155203?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880635234155342?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880634607488738?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880519440833588?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880251330860399?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880249834193882?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880249220860610?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880135807538618?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/880018180883714?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879993847552814?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879975050888027?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879828897569309?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879825997569599?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879826507569548?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879494384269427?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879491310936401?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879227444296121?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879226637629535?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879225547629644?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879224727629726?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879223810963151?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879102794308586?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879101247642074?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/879099590975573?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/878804271005105?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/878430701042462?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/878554227696776?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/878327984386067?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877943121091220?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877839401101592?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877837871101745?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877593647792834?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877592734459592?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877591714459694?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877590381126494?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/877588994459966?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876932711192261?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876934391192093?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876933934525472?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876933481192184?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876931847859014?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876842871201245?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876842197867979?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876413567910842?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876841547868044?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876236571261875?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876235047928694?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876234261262106?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876233224595543?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876024484616417?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/876022357949963?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/875568004662065?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/875568791328653?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/875359594682906?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/875360158016183?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/875359028016296?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874848064734059?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874847301400802?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874846491400883?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874593441426188?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874422351443297?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874421478110051?ref=embed_post",
    "https://www.facebook.com/thoisuvtv/posts/874420481443484?ref=embed_post"
)

$startRow = 25
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $urls[$i]
}

$lastRow = $startRow + $urls.Length - 1

$ws.Application.ActiveWindow.ScrollRow = 94
$ws.Range("A$lastRow").Select()
